# Updated symbol list (crypto price/volume/hour snapshot refresh)
# Applies the per-cell value updates described by the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price), E (Volume %) and G (Hora) hold numeric-looking text
# (e.g. "247.01", "1.43%", "9") that must stay stored as plain text, exactly
# like the rest of the sheet. Pre-formatting these cells as Text ("@") before
# writing the values stops Excel from auto-converting them into real numbers
# or percentages.
$textForceCells = @('D2', 'E2', 'G2', 'D3', 'E3', 'G3', 'D4', 'E4', 'G4', 'D5', 'E5', 'G5', 'D6', 'E6', 'G6', 'D7', 'E7', 'G7', 'D8', 'E8', 'G8', 'E9', 'G9', 'D10', 'E10', 'G10', 'D11', 'E11', 'G11', 'D12', 'E12', 'G12', 'D13', 'E13', 'G13', 'D14', 'E14', 'G14', 'D15', 'E15', 'G15', 'D16', 'E16', 'G16', 'D17', 'E17', 'G17', 'D18', 'E18', 'G18', 'D19', 'E19', 'G19', 'D20', 'E20', 'G20', 'D21', 'E21', 'G21', 'E22', 'G22', 'D23', 'E23', 'G23', 'E24', 'G24', 'D25', 'E25', 'G25', 'D26', 'E26', 'G26', 'D27', 'E27', 'G27', 'G28', 'G29', 'G30', 'G31', 'G32', 'G33', 'G34', 'G35', 'G36', 'G37', 'G38', 'G39', 'D40', 'E40', 'G40', 'D41', 'E41', 'G41', 'D42', 'E42', 'G42', 'D43', 'E43', 'G43', 'D44', 'E44', 'G44', 'D45', 'E45', 'G45', 'E46', 'G46', 'G47', 'D48', 'E48', 'G48', 'E49', 'G49', 'E50', 'G50', 'G51')
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Updated values (row-by-row, left to right) ---
# Row 2
$ws.Range('D2').Value = '247.01'
$ws.Range('E2').Value = '1.43%'
$ws.Range('G2').Value = '9'
# Row 3
$ws.Range('D3').Value = '30.15'
$ws.Range('E3').Value = '12.01%'
$ws.Range('G3').Value = '9'
# Row 4
$ws.Range('D4').Value = '5.180'
$ws.Range('E4').Value = '0.47%'
$ws.Range('G4').Value = '9'
# Row 5
$ws.Range('D5').Value = '0.05750'
$ws.Range('E5').Value = '2.34%'
$ws.Range('G5').Value = '9'
# Row 6
$ws.Range('D6').Value = '6.587'
$ws.Range('E6').Value = '1.49%'
$ws.Range('G6').Value = '9'
# Row 7
$ws.Range('D7').Value = '0.8579'
$ws.Range('E7').Value = '4.96%'
$ws.Range('G7').Value = '9'
# Row 8
$ws.Range('D8').Value = '0.8805'
$ws.Range('E8').Value = '6.08%'
$ws.Range('G8').Value = '9'
# Row 9
$ws.Range('E9').Value = '3.02%'
$ws.Range('G9').Value = '9'
# Row 10
$ws.Range('D10').Value = '0.07005'
$ws.Range('E10').Value = '1.52%'
$ws.Range('G10').Value = '9'
# Row 11
$ws.Range('D11').Value = '0.02915'
$ws.Range('E11').Value = '0.67%'
$ws.Range('G11').Value = '9'
# Row 12
$ws.Range('D12').Value = '0.09390'
$ws.Range('E12').Value = '0.16%'
$ws.Range('G12').Value = '9'
# Row 13
$ws.Range('D13').Value = '0.001521'
$ws.Range('E13').Value = '0.56%'
$ws.Range('G13').Value = '9'
# Row 14
$ws.Range('D14').Value = '0.04123'
$ws.Range('E14').Value = '-10.18%'
$ws.Range('G14').Value = '9'
# Row 15
$ws.Range('B15').Value = 'One'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D15').Value = '0.0006012'
$ws.Range('E15').Value = '0.53%'
$ws.Range('G15').Value = '9'
# Row 16
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '0.005988'
$ws.Range('E16').Value = '-2.48%'
$ws.Range('G16').Value = '9'
# Row 17
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.504'
$ws.Range('E17').Value = '-2.95%'
$ws.Range('G17').Value = '9'
# Row 18
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').Value = '3.081'
$ws.Range('E18').Value = '1.87%'
$ws.Range('G18').Value = '9'
# Row 19
$ws.Range('B19').Value = 'BTSEToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D19').Value = '2.174'
$ws.Range('E19').Value = '-2.27%'
$ws.Range('G19').Value = '9'
# Row 20
$ws.Range('D20').Value = '0.3145'
$ws.Range('E20').Value = '1.06%'
$ws.Range('G20').Value = '9'
# Row 21
$ws.Range('D21').Value = '0.03296'
$ws.Range('E21').Value = '6.84%'
$ws.Range('G21').Value = '9'
# Row 22
$ws.Range('E22').Value = '1.16%'
$ws.Range('G22').Value = '9'
# Row 23
$ws.Range('D23').Value = '3.600'
$ws.Range('E23').Value = '-4.15%'
$ws.Range('G23').Value = '9'
# Row 24
$ws.Range('E24').Value = '2.71%'
$ws.Range('G24').Value = '9'
# Row 25
$ws.Range('D25').Value = '0.001212'
$ws.Range('E25').Value = '-1.03%'
$ws.Range('G25').Value = '9'
# Row 26
$ws.Range('D26').Value = '0.004505'
$ws.Range('E26').Value = '0.42%'
$ws.Range('G26').Value = '9'
# Row 27
$ws.Range('D27').Value = '0.0001179'
$ws.Range('E27').Value = '20.28%'
$ws.Range('G27').Value = '9'
# Row 28
$ws.Range('G28').Value = '9'
# Row 29
$ws.Range('G29').Value = '9'
# Row 30
$ws.Range('G30').Value = '9'
# Row 31
$ws.Range('G31').Value = '9'
# Row 32
$ws.Range('G32').Value = '9'
# Row 33
$ws.Range('G33').Value = '9'
# Row 34
$ws.Range('G34').Value = '9'
# Row 35
$ws.Range('G35').Value = '9'
# Row 36
$ws.Range('G36').Value = '9'
# Row 37
$ws.Range('G37').Value = '9'
# Row 38
$ws.Range('G38').Value = '9'
# Row 39
$ws.Range('G39').Value = '9'
# Row 40
$ws.Range('D40').Value = '0.03790'
$ws.Range('E40').Value = '4.20%'
$ws.Range('G40').Value = '9'
# Row 41
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').Value = '0.1070'
$ws.Range('E41').Value = '-22.25%'
$ws.Range('G41').Value = '9'
# Row 42
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').Value = '0.002587'
$ws.Range('E42').Value = '0.28%'
$ws.Range('G42').Value = '9'
# Row 43
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D43').Value = '0.003516'
$ws.Range('E43').Value = '-41.98%'
$ws.Range('G43').Value = '9'
# Row 44
$ws.Range('D44').Value = '0.01008'
$ws.Range('E44').Value = '21.07%'
$ws.Range('G44').Value = '9'
# Row 45
$ws.Range('D45').Value = '0.00005111'
$ws.Range('E45').Value = '-3.80%'
$ws.Range('G45').Value = '9'
# Row 46
$ws.Range('E46').Value = '-0.14%'
$ws.Range('G46').Value = '9'
# Row 47
$ws.Range('G47').Value = '9'
# Row 48
$ws.Range('D48').Value = '0.002718'
$ws.Range('E48').Value = '3.85%'
$ws.Range('G48').Value = '9'
# Row 49
$ws.Range('E49').Value = '-0.14%'
$ws.Range('G49').Value = '9'
# Row 50
$ws.Range('E50').Value = '-0.14%'
$ws.Range('G50').Value = '9'
# Row 51
$ws.Range('G51').Value = '9'
